$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (Q1, R1) with the same header style as the existing header row ---
$ws.Range("Q1").Value = "EENS 95% CI"
$ws.Range("R1").Value = "EENS 99% CI"
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1:R1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 2 (LP1 / residential) ---
$ws.Range("F2").Value = 2.097658007960921
$ws.Range("G2").Value = 6608
$ws.Range("H2").Value = 0.9701033402434285
$ws.Range("I2").Value = 2.162303664921466
$ws.Range("J2").Value = 454.0837696335079
$ws.Range("K2").Value = 440.5081816717935
$ws.Range("L2").Value = 203.72170145112
$ws.Range("M2").Value = 1.122247034259093

# --- Row 3 (LP2 / residential) ---
$ws.Range("F3").Value = 3.046956170178304
$ws.Range("G3").Value = 6608
$ws.Range("H3").Value = 1.409125008484397
$ws.Range("I3").Value = 2.162303664921466
$ws.Range("J3").Value = 454.0837696335079
$ws.Range("K3").Value = 639.8607957374438
$ws.Range("L3").Value = 295.9162517817234
$ws.Range("M3").Value = 1.630121551045393

# --- Row 4 (LP3 / residential) ---
$ws.Range("F4").Value = 3.853935416931458
$ws.Range("G4").Value = 6608
$ws.Range("H4").Value = 1.782328485796389
$ws.Range("I4").Value = 2.162303664921466
$ws.Range("J4").Value = 454.0837696335079
$ws.Range("K4").Value = 809.3264375556062
$ws.Range("L4").Value = 374.2889820172416
$ws.Range("M4").Value = 2.06185544805833

# --- Row 5 (LP4, now industrial instead of residential) ---
$ws.Range("B5").Value = "industrial"
$ws.Range("F5").Value = 4.221194677684629
$ws.Range("G5").Value = 6608
$ws.Range("H5").Value = 1.952174778299671
$ws.Range("I5").Value = 2.162303664921466
$ws.Range("J5").Value = 2.162303664921466
$ws.Range("K5").Value = 4.221194677684629
$ws.Range("L5").Value = 1.952174778299671
$ws.Range("M5").Value = 2.3891961875695

# --- Row 6 (TOTAL) ---
$ws.Range("J6").Value = 2.162303664921466
$ws.Range("K6").Value = 3.001452630178333
$ws.Range("L6").Value = 1.388080998460198
$ws.Range("M6").Value = 7.203420220932315
$ws.Range("N6").Value = 3056
$ws.Range("P6").Value = 0.01993630671266165
$ws.Range("Q6").Value = "(6.921994197234212, 7.484947776640266)"
$ws.Range("R6").Value = "(6.833530063327546, 7.573411910546931)"
